$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.037.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.439.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.737"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000221"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.984.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.435.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.151.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "495.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +22.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "93.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "34.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.39%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.12%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.74%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.29%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.138"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.320"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +24.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.69%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +24.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.148"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.69%  "
